# "break out stock.yaml completed"
# Target sheet: "10per change" (first sheet in the workbook)
#  - E48 was mistakenly stored as text "590024"; fix it up to the numeric value 590024
#    (every other bsecode in the column is numeric)
#  - Append a newly-scraped row 49 for the same FACT record (second scrape pass that day),
#    whose bsecode came through as text "590024" just like E48 did before the fix

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("10per change")

# Fix E48: was text, should be numeric 590024
$ws.Cells.Item(48, 5).Value = 590024

# New row 49
$ws.Cells.Item(49, 1).Value = "25/06/2024 05:45:46"
$ws.Cells.Item(49, 2).Value = 1
$ws.Cells.Item(49, 3).Value = "FACT"
$ws.Cells.Item(49, 4).Value = "Fertilizers And Chemicals Travancore Limited"

# Force E49 to stay a text value "590024" (matches the source scrape) instead of
# being auto-coerced to a number, then drop back to the default style so no
# stray number-format style is left attached to the cell.
$ws.Cells.Item(49, 5).NumberFormat = "@"
$ws.Cells.Item(49, 5).Value = "590024"
$ws.Cells.Item(49, 5).Style = "Normal"

$ws.Cells.Item(49, 6).Value = -1.04
$ws.Cells.Item(49, 7).Value = 1010.35
$ws.Cells.Item(49, 8).Value = 1217238
